# "remove - innecesary stuff"
# Adds a new "Sesion por tema" summary sheet (Sheet2) after Sheet1, moves the
# active-sheet/selection focus to it, and leaves Sheet1's own selection at M13.

$wb = $excel.ActiveWorkbook

# --- Sheet1: drop its "tab selected" state by moving its own cell cursor
# first (still the active sheet at this point), before Sheet2 steals focus.
$ws1 = $wb.Worksheets.Item(1)
$ws1.Range("M13").Select()

# --- Add the new sheet right after Sheet1 ---
$ws2 = $wb.Worksheets.Add($null, $ws1)

# --- Header row ---
$ws2.Range("A1").Value = "Temas"
$ws2.Range("B1").Value = "Sesion"
$ws2.Range("C1").Value = "Estimado"

# --- Topic column ---
$ws2.Range("A2").Value = "Conceptos Basicos"
$ws2.Range("A3").Value = "Introduccion a Jmeter"
$ws2.Range("A4").Value = "Instalacion de Jmeter"
$ws2.Range("A5").Value = "Plugins Manager"
$ws2.Range("A6").Value = "Recording "
$ws2.Range("A7").Value = "Configuración de CSV Data set"
$ws2.Range("A8").Value = "Expresiones regulares"
$ws2.Range("A9").Value = "Jmeter funciones"
$ws2.Range("A10").Value = "Ejecucion de pruebas"
$ws2.Range("A11").Value = "Analisis de resultados"

# --- Session / estimate columns ---
$ws2.Range("B2").Value = "S1"
$ws2.Range("C2").Value = "4h"
$ws2.Range("B5").Value = "S2"
$ws2.Range("C5").Value = "4h"

# --- Column width (≈ width 28 in the saved sheet) ---
$ws2.Columns.Item(1).ColumnWidth = 27.15

# --- Alignment: Estimado column is centered + wrapped ---
$r1 = $ws2.Range("C2:C4")
$r1.HorizontalAlignment = -4108
$r1.VerticalAlignment = -4108
$r1.WrapText = $true

# --- Alignment: Sesion column (and the second Estimado block) centered, no wrap ---
$r2 = $ws2.Range("B2:B4")
$r2.HorizontalAlignment = -4108
$r2.VerticalAlignment = -4108

$r3 = $ws2.Range("B5:B6")
$r3.HorizontalAlignment = -4108
$r3.VerticalAlignment = -4108

$r4 = $ws2.Range("C5:C6")
$r4.HorizontalAlignment = -4108
$r4.VerticalAlignment = -4108

# --- Merge the session/estimate blocks ---
$ws2.Range("B2:B4").Merge()
$ws2.Range("C2:C4").Merge()
$ws2.Range("B5:B6").Merge()
$ws2.Range("C5:C6").Merge()

# --- Selection that ends up saved for the new sheet ---
$ws2.Range("B7").Select()
